# Generate Report for Handoff
#
# Rows whose "Ready for handoff" status was (re)generated during this
# handoff pass are rows 7, 9, 10, 11, 12, 14 on each localization sheet.
# For those rows we:
#   - bump the handoff timestamps forward (report regenerated a little later)
#   - stamp the "Priority" column with "ht" (matches row 6, the other row
#     that already went through a handoff in this file)

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-20 12:17:45"
}

# --- zh-cn / de-de sheets: "Priority" (column E) + "Latest Handoff Datetime" (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-20 12:17:40"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    # de-de's "Latest Handoff Datetime" shares its text with Overview's
    # "Latest HO Xliff Generate Date" column, so it takes the same stamp.
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-20 12:17:45"
}
